$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 13, mirroring row 12 (CIMS / Root / CIMS / Service requested / n/a / unit)
# but with Target = "CIMS.RoW" (a new region alongside existing "CIMS.CAN").
$ws.Range("A13").Value = "CIMS"
$ws.Range("B13").Value = "Root"
$ws.Range("C13").Value = "CIMS"
$ws.Range("G13").Value = "Service requested"
$ws.Range("J13").Value = "CIMS.RoW"
$ws.Range("K13").Value = "n/a"
$ws.Range("L13").Value = "unit"

$ws.Range("M13:W13").Value = 1

# Update selection to match target state
$ws.Range("J14").Select()
